$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 - MIT Environmental Solutions Initative (John Fernandez)
$ws.Range("B27").Value = "None,AIR-INK: Air-Pollution to ink"

# Row 38 - The Pershing Square Foundation
$ws.Range("B38").Value = "None,Renewal Workshop"

# Row 27 continued
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "this is it!"

# Row 38 continued
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = "Maybe not"
